$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 4207
$ws.Range("I21").Value = 4658.1665
$ws.Range("J21").Value = 1500
$ws.Range("K21").Value = 4658.1665
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = -4190.1665
$ws.Range("N21").Value = -2436
$ws.Range("H23").Value = 4207
$ws.Range("I23").Value = 4658.1665
$ws.Range("J23").Value = 1500
$ws.Range("K23").Value = 4658.1665
$ws.Range("L23").Value = 1500
$ws.Range("M23").Value = -4424.1665
$ws.Range("N23").Value = -1968
$ws.Range("H64").Value = 4875
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null
$ws.Range("H67").Value = 4875
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null
$ws.Range("H107").Value = 4226.923
$ws.Range("I107").Value = 4062.8262
$ws.Range("J107").Value = 5485
$ws.Range("K107").Value = 4062.8262
$ws.Range("L107").Value = 5485
$ws.Range("M107").Value = -2142.8262
$ws.Range("N107").Value = -9325
$ws.Range("H111").Value = 869.3333
$ws.Range("I111").Value = 869.3333
$ws.Range("K111").Value = 2607.9999
$ws.Range("M111").Value = 459.0001000000002
$ws.Range("H113").Value = 2775.3
$ws.Range("I113").Value = 2112.8
$ws.Range("J113").Value = 3437.8
$ws.Range("K113").Value = 2112.8
$ws.Range("L113").Value = 3437.8
$ws.Range("M113").Value = 1141.2
$ws.Range("N113").Value = -9945.799999999999
$ws.Range("H115").Value = 1360.2
$ws.Range("I115").Value = 887.75
$ws.Range("J115").Value = 3250
$ws.Range("K115").Value = 2663.25
$ws.Range("L115").Value = 9750
$ws.Range("M115").Value = -1096.25
$ws.Range("N115").Value = -12884
$ws.Range("H138").Value = 2240.5845
$ws.Range("J138").Value = 2489.7693
$ws.Range("L138").Value = 7469.3079
$ws.Range("N138").Value = -17749.3079
$ws.Range("H141").Value = 4697.909
$ws.Range("I141").Value = 5158
$ws.Range("J141").Value = 4314.5
$ws.Range("K141").Value = 15474
$ws.Range("L141").Value = 12943.5
$ws.Range("M141").Value = -10294
$ws.Range("N141").Value = -23303.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1850
$ws.Range("H33").Value = 5675.3335
$ws.Range("J33").Value = 6000
$ws.Range("L33").Value = 6000
$ws.Range("N33").Value = -6658
$ws.Range("H102").Value = 4121.1665
$ws.Range("I102").Value = 3404.9092
$ws.Range("K102").Value = 3404.9092
$ws.Range("M102").Value = -1782.9092
$ws.Range("H116").Value = 1850
$ws.Range("H132").Value = 2780466.2
$ws.Range("I132").Value = 3127649.5
$ws.Range("K132").Value = 9382948.5
$ws.Range("M132").Value = -9380418.5
$ws.Range("H134").Value = 64999
$ws.Range("J134").Value = 64999
$ws.Range("L134").Value = 64999
$ws.Range("N134").Value = -75139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1850
$ws.Range("H19").Value = 990
$ws.Range("J19").Value = 1066.6666
$ws.Range("L19").Value = 1066.6666
$ws.Range("N19").Value = -1412.6666
$ws.Range("H80").Value = 1946
$ws.Range("I80").Value = 1543
$ws.Range("K80").Value = 1543
$ws.Range("M80").Value = -545
$ws.Range("H83").Value = 1946
$ws.Range("I83").Value = 1543
$ws.Range("K83").Value = 7715
$ws.Range("M83").Value = -2723
$ws.Range("H86").Value = 4400.25
$ws.Range("I86").Value = 4311.4443
$ws.Range("K86").Value = 4311.4443
$ws.Range("M86").Value = -3188.4443
$ws.Range("H89").Value = 4400.25
$ws.Range("I89").Value = 4311.4443
$ws.Range("K89").Value = 21557.2215
$ws.Range("M89").Value = -15941.2215
$ws.Range("H94").Value = 1779.5
$ws.Range("I94").Value = 1636.561
$ws.Range("K94").Value = 1636.561
$ws.Range("M94").Value = -1185.561
$ws.Range("H99").Value = 9171.154
$ws.Range("I99").Value = 12507.889
$ws.Range("K99").Value = 12507.889
$ws.Range("M99").Value = -11009.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1999.6666
$ws.Range("I16").Value = 1999.6666
$ws.Range("K16").Value = 1999.6666
$ws.Range("M16").Value = -1712.6666
$ws.Range("H31").Value = 2487.182
$ws.Range("I31").Value = 2278.2632
$ws.Range("J31").Value = 2770.7144
$ws.Range("K31").Value = 2278.2632
$ws.Range("L31").Value = 2770.7144
$ws.Range("M31").Value = -1983.2632
$ws.Range("N31").Value = -3360.7144
$ws.Range("H34").Value = 2487.182
$ws.Range("I34").Value = 2278.2632
$ws.Range("J34").Value = 2770.7144
$ws.Range("K34").Value = 2278.2632
$ws.Range("L34").Value = 2770.7144
$ws.Range("M34").Value = -2076.2632
$ws.Range("N34").Value = -3174.7144
$ws.Range("H62").Value = 4571.4287
$ws.Range("I62").Value = 4437.8
$ws.Range("K62").Value = 4437.8
$ws.Range("M62").Value = -3813.8
$ws.Range("H65").Value = 4571.4287
$ws.Range("I65").Value = 4437.8
$ws.Range("K65").Value = 22189
$ws.Range("M65").Value = -19069
$ws.Range("H113").Value = 1999.6666
$ws.Range("I113").Value = 1999.6666
$ws.Range("K113").Value = 1999.6666
$ws.Range("M113").Value = 170.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5834870
$ws.Range("I4").Value = 8750259
$ws.Range("J4").Value = 4093.5
$ws.Range("K4").Value = 26250777
$ws.Range("L4").Value = 12280.5
$ws.Range("M4").Value = -26250665
$ws.Range("N4").Value = -12504.5
$ws.Range("H107").Value = 3284.3333
$ws.Range("J107").Value = 4067.375
$ws.Range("L107").Value = 12202.125
$ws.Range("N107").Value = -16042.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1579.5
$ws.Range("I80").Value = 1465
$ws.Range("J80").Value = 1665.375
$ws.Range("K80").Value = 1465
$ws.Range("L80").Value = 1665.375
$ws.Range("M80").Value = -467
$ws.Range("N80").Value = -3661.375
$ws.Range("H83").Value = 1579.5
$ws.Range("I83").Value = 1465
$ws.Range("J83").Value = 1665.375
$ws.Range("K83").Value = 7325
$ws.Range("L83").Value = 8326.875
$ws.Range("M83").Value = -2333
$ws.Range("N83").Value = -18310.875
$ws.Range("H107").Value = 5696.3335
$ws.Range("I107").Value = 9709.637000000001
$ws.Range("J107").Value = 1281.7
$ws.Range("K107").Value = 9709.637000000001
$ws.Range("L107").Value = 1281.7
$ws.Range("M107").Value = -7789.637000000001
$ws.Range("N107").Value = -5121.7
$ws.Range("H122").Value = 63043.47
$ws.Range("I122").Value = 103874.2
$ws.Range("J122").Value = 4713.857
$ws.Range("K122").Value = 311622.6
$ws.Range("L122").Value = 14141.571
$ws.Range("M122").Value = -309172.6
$ws.Range("N122").Value = -19041.571
$ws.Range("H132").Value = 7278.75
$ws.Range("I132").Value = 8828.546
$ws.Range("K132").Value = 26485.638
$ws.Range("M132").Value = -23955.638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6161.857
$ws.Range("I7").Value = 4626.8
$ws.Range("K7").Value = 4626.8
$ws.Range("M7").Value = -4514.8
$ws.Range("H16").Value = 395.5
$ws.Range("I16").Value = 395.5
$ws.Range("K16").Value = 395.5
$ws.Range("M16").Value = -225.5
$ws.Range("H40").Value = 4339.25
$ws.Range("I40").Value = 4207.1
$ws.Range("K40").Value = 4207.1
$ws.Range("M40").Value = -4071.1
$ws.Range("H46").Value = 6556.7144
$ws.Range("J46").Value = 7499.8335
$ws.Range("L46").Value = 7499.8335
$ws.Range("N46").Value = -7875.8335
$ws.Range("H93").Value = 2965.2307
$ws.Range("I93").Value = 1300.6666
$ws.Range("J93").Value = 3464.6
$ws.Range("K93").Value = 1300.6666
$ws.Range("L93").Value = 3464.6
$ws.Range("M93").Value = -52.66660000000002
$ws.Range("N93").Value = -5960.6
$ws.Range("H126").Value = 6161.857
$ws.Range("I126").Value = 4626.8
$ws.Range("K126").Value = 13880.4
$ws.Range("M126").Value = -11410.4
$ws.Range("H136").Value = 19896
$ws.Range("I136").Value = 6854.4
$ws.Range("K136").Value = 20563.2
$ws.Range("M136").Value = -18013.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 20736
$ws.Range("I32").Value = 14981.333
$ws.Range("K32").Value = 14981.333
$ws.Range("M32").Value = -14664.333
$ws.Range("H34").Value = 29
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 29
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 29
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -435
$ws.Range("H64").Value = 49000
$ws.Range("J64").Value = 49000
$ws.Range("L64").Value = 49000
$ws.Range("N64").Value = -49496
$ws.Range("H67").Value = 49000
$ws.Range("J67").Value = 49000
$ws.Range("L67").Value = 49000
$ws.Range("N67").Value = -50716
$ws.Range("H113").Value = 940.4
$ws.Range("J113").Value = 887
$ws.Range("L113").Value = 2661
$ws.Range("N113").Value = -7001
$ws.Range("H122").Value = 114235.6
$ws.Range("I122").Value = 1880
$ws.Range("K122").Value = 5640
$ws.Range("M122").Value = -3190
$ws.Range("H132").Value = 3070.1875
$ws.Range("I132").Value = 3076.9167
$ws.Range("J132").Value = 3050
$ws.Range("K132").Value = 9230.750100000001
$ws.Range("L132").Value = 9150
$ws.Range("M132").Value = -6700.750100000001
$ws.Range("N132").Value = -14210
$ws.Range("H136").Value = 3900.375
$ws.Range("I136").Value = 3743.2856
$ws.Range("K136").Value = 11229.8568
$ws.Range("M136").Value = -8679.856800000001
